# Update Car_Circuitry_BOM workbook per "update bom & things" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Car_Circuitry_BOM")

# Helper: set a cell's text value while preserving its existing cell
# format (setting .Value alone can reset the "quote prefix" flag that
# is baked into this sheet's cell style).
function Set-TextPreserveFormat($ws, $targetAddr, $formatSourceAddr, $text) {
    $target = $ws.Range($targetAddr)
    $target.Value = $text
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
}

# Row 23: Designator "R3, R6" -> "R3, R6, R10"; Quantity 2 -> 3
Set-TextPreserveFormat $ws "C23" "B23" "R3, R6, R10"
$ws.Range("F23").Value = 3

# Row 25: Comment (blank) -> "301 Ω"
Set-TextPreserveFormat $ws "A25" "B25" "301 Ω"

# Row 27: Designator "R10" -> "R11"
Set-TextPreserveFormat $ws "C27" "B27" "R11"

# Row 28: Designator "R11" -> "R12"
Set-TextPreserveFormat $ws "C28" "B28" "R12"

$excel.CutCopyMode = 0
